$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "description" column before the existing "reason" column (C),
# so the old C-column (description text) shifts to C and "reason" becomes D.
$ws.Range("C1").EntireColumn.Insert()

# Header row: fill new C1 and relabel the shifted D1.
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "reason"

# Insert a new row above the current last data row (old row 4, "id" 3 /
# NLP engineer) so there is room for the extra "id" 2 record, giving a
# total of four data rows (rows 2-5).
$ws.Range("A4").EntireRow.Insert()

# Column A ("id") holds text values, not numbers - force text storage so
# values like "4" aren't auto-coerced to numeric.
$ws.Range("A2:A5").NumberFormat = "@"

# Row 2 (id=4, score=100)
$ws.Cells.Item(2, 1).Value = "4"
$ws.Cells.Item(2, 2).Value = 100
$ws.Cells.Item(2, 3).Value = "Application Development: Front-end coding for a sweat amino acid analysis app under a Scrum Master. Skills: Flutter, Dart, Android/iOS, Firestore, Firebase Authentication, Cloud Storage/Messaging, Mobile app architecture/design, Git."
$ws.Cells.Item(2, 4).Value = "The job description requires front-end coding for a sweat amino acid analysis app using Flutter, Dart, and various Firebase services. Your experience in developing a food delivery platform using Flutter, Dart, and Firebase demonstrates your proficiency in the required skills. Therefore, the job is suitable for you."

# Row 3 (id=1, score=88)
$ws.Cells.Item(3, 1).Value = "1"
$ws.Cells.Item(3, 2).Value = 88
$ws.Cells.Item(3, 3).Value = "SDE Intern: Remote MERN Stack Developer Internship, responsible for designing and developing web/mobile applications using MongoDB, ExpressJS, ReactJS, and NodeJS. Additional tasks include code maintenance, scalability, feature development, and product enhancement suggestions. Skills required: MongoDB, ReactJS, JavaScript, Web Development, NodeJS."
$ws.Cells.Item(3, 4).Value = "The job as a Remote MERN Stack Developer Internship matches well with your experience in developing web/mobile applications using MongoDB, ExpressJS, ReactJS, and NodeJS. Your skills in MongoDB, ReactJS, JavaScript, and Web Development make you well suited for the tasks of code maintenance, scalability, feature development, and product enhancement."

# Row 4 (id=2, score=69)
$ws.Cells.Item(4, 1).Value = "2"
$ws.Cells.Item(4, 2).Value = 69
$ws.Cells.Item(4, 3).Value = "Frontend Engineer Intern - Work in a team to ensure consistent web design and user experience, optimize web pages, and maintain brand consistency. Requires excellent communication skills and proficiency in ReactJS, JavaScript, CSS, and NextJS. 3-month evaluative unpaid internship with potential return offers."
$ws.Cells.Item(4, 4).Value = "The job as a Frontend Engineer Intern at the company is suitable for you because it requires excellent communication skills and proficiency in ReactJS, JavaScript, CSS, and NextJS, which align with your experience in developing web applications using these technologies. The evaluative unpaid internship also provides an opportunity to gain valuable experience and potentially secure return offers. The moderate score suggests that while it may not be a perfect fit, it is still a suitable opportunity for you."

# Row 5 (id=3, score=30.5)
$ws.Cells.Item(5, 1).Value = "3"
$ws.Cells.Item(5, 2).Value = 30.5
$ws.Cells.Item(5, 3).Value = "NLP Engineer (Remote): Collect and preprocess text corpora for language model training. Analyze data, develop and improve models. Skills: NLP, Pytorch, Computer Vision, Python."
$ws.Cells.Item(5, 4).Value = "Based on your projects and the job description, the role of NLP Engineer (Remote) seems moderately suitable for you. While you have experience in utilizing Python and have knowledge of NLP and Pytorch, your projects do not directly align with the job requirements of collecting and preprocessing text corpora for language model training."
